$p = $ppt.ActivePresentation

# Delete the last slide (slide id 280, slide5.xml) - it was a wrongly named
# duplicate figure that doesn't belong in the deck.
$p.Slides.Item(5).Delete()

# Clean up the "byte 0" and "byte 2" labels: their paragraphs carried a
# redundant trailing endParaRPr (identical formatting to the run itself).
# Deleting the text range first and retyping the text drops that stray
# endParaRPr, matching the other "byte N" labels on the slide.
$s4 = $p.Slides.Item(4)

$tr1 = $s4.Shapes.Item("Rectangle 6").TextFrame.TextRange
$tr1.Delete()
$tr1.Text = "byte 0"

$tr2 = $s4.Shapes.Item("Rectangle 12").TextFrame.TextRange
$tr2.Delete()
$tr2.Text = "byte 2"
